$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 39, shifting the existing rows 39:135 down to 40:136
$ws.Rows("39").Insert()

# Populate the newly inserted row 39 with the new weekly price entry.
# A:C keep the same market/region info used throughout the sheet.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = Get-Date -Year 2022 -Month 3 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100114007
$ws.Range("G39").Value = "Jengibre"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 40
$ws.Range("K39").Value = 26000
$ws.Range("L39").Value = 26000
$ws.Range("M39").Value = 26000
$ws.Range("N39").Value = "`$/caja 13 kilos"
$ws.Range("O39").Value = "Perú"
$ws.Range("P39").Value = 2000
$ws.Range("Q39").Value = 13
$ws.Range("R39").Value = "Hortaliza"
